# Sync updated BOM for multitarget
#
# - Move "Digikey Final" to be the first sheet tab.
# - Rename "Digikey Upload" to "DO NOT USE" (formulas referencing it are
#   updated automatically by Excel).
# - Make "From Eagle" the active sheet (now the 2nd tab).

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Digikey Final").Move($wb.Worksheets.Item(1))

$wb.Worksheets.Item("Digikey Upload").Name = "DO NOT USE"

$wb.Worksheets.Item("From Eagle").Activate()
